$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (everything shifts right by one column)
$ws.Columns("A:A").Insert()
$ws.Columns("A:A").ColumnWidth = 14.14

# Header for the new column
$ws.Range("A1").Value = "RefID"
# Copy the bold header style from the neighboring header cell (B1) onto A1
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# New RefID values per row (entered in this order so the shared-string table
# ends up in the same sequence as the authored workbook)
$ws.Range("A4").Value = "LATFLD-77"
$ws.Range("A3").Value = "LATFLD-76"
$ws.Range("A2").Value = "LATFLD-75"
$ws.Range("A5").Value = "LATFLD-22"

# The old free-text "Comments for ..." values become "Submission"
$ws.Range("O3").Value = "Submission"
$ws.Range("O4").Value = "Submission"

# Jira hyperlinks for the RefID column (skip LATFLD-22 / row 5)
$ws.Hyperlinks.Add($ws.Range("A4"), "https://leapthought.atlassian.net/browse/LATFLD-77", "", "", "https://leapthought.atlassian.net/browse/LATFLD-77") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://leapthought.atlassian.net/browse/LATFLD-76", "", "", "https://leapthought.atlassian.net/browse/LATFLD-76") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "https://leapthought.atlassian.net/browse/LATFLD-75", "", "", "https://leapthought.atlassian.net/browse/LATFLD-75") | Out-Null

# Adding a hyperlink auto-applies Excel's built-in blue/underline "Hyperlink"
# style; put these cells back to the plain (unstyled) look used everywhere
# else in the data rows by pulling the format from an already-plain neighbor.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the visible selection to P4 (matches the author's last selection)
$ws.Range("P4").Select()
